$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pistures")
$ws.Range("A13").Value = "Ground"
$ws.Range("E13").Value = "http://creativecommons.org/publicdomain/mark/1.0/"
$ws.Range("D13").Value = "Public Domain Mark 1.0"
$ws.Range("C13").Value = "No Copyright"
$ws.Range("F13").Value = "https://www.makeschool.com/academy/art/object/assorted-clouds"
$ws.Range("F13").Select()
